$d = $word.ActiveDocument

# Locate the paragraph "Define the term eutrophication." via Find.
$findRange = $d.Content
$found = $findRange.Find.Execute("Define the term eutrophication.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target paragraph 'Define the term eutrophication.'"
}

$targetParagraph = $findRange.Paragraphs.Item(1)

# Determine the paragraph's 1-based index within the document's Paragraphs collection.
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $targetParagraph.Range.Start) {
        $targetIndex = $i
        break
    }
}
if ($targetIndex -eq 0) {
    throw "Could not resolve paragraph index for target paragraph"
}

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# A paragraph containing only a manual page break.
$pageBreakXml = "<w:p $wNs>" + `
    "<w:pPr>" + `
        "<w:spacing w:line='259' w:lineRule='auto'/>" + `
        "<w:jc w:val='left'/>" + `
        "<w:rPr><w:rFonts w:cs='Times New Roman'/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr>" + `
    "</w:pPr>" + `
    "<w:r>" + `
        "<w:rPr><w:rFonts w:cs='Times New Roman'/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr>" + `
        "<w:br w:type='page'/>" + `
    "</w:r>" + `
    "</w:p>"

# An empty List Paragraph (matches the pre-existing blank line already in the doc).
$emptyListParaXml = "<w:p $wNs>" + `
    "<w:pPr>" + `
        "<w:pStyle w:val='ListParagraph'/>" + `
        "<w:ind w:left='1530'/>" + `
        "<w:rPr><w:rFonts w:cs='Times New Roman'/><w:sz w:val='28'/><w:szCs w:val='28'/></w:rPr>" + `
    "</w:pPr>" + `
    "</w:p>"

# Insert a fresh paragraph right after the target paragraph, then stamp it
# with the page-break XML.
$r = $d.Paragraphs.Item($targetIndex).Range
$r.Collapse(0)
$newPara1 = $r.InsertParagraphAfter()
$discard = $d.Paragraphs.Item($targetIndex + 1).Range.InsertXML($pageBreakXml)

# Insert four empty ListParagraph paragraphs after the page break, each one
# right after the previous newly-created paragraph.
for ($j = 1; $j -le 4; $j++) {
    $rr = $d.Paragraphs.Item($targetIndex + $j).Range
    $rr.Collapse(0)
    $newParaN = $rr.InsertParagraphAfter()
    $discard = $d.Paragraphs.Item($targetIndex + $j + 1).Range.InsertXML($emptyListParaXml)
}

Write-Output "Inserted page break + 4 blank paragraphs after paragraph $targetIndex. New paragraph count: $($d.Paragraphs.Count)"
